$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-export of the speaker-variant playlist with no "is_prefered" (D) flag
# and no Levenshtein-distance-based id reuse: ids are now 1:1 with their
# original slug (duplicated slugs in the old export, e.g. repeated
# "#1.-vrou"/"#1.-borg"/"#1.-sold", are replaced by their own distinct ids)
# and the rows come out in a different order.
$rows = @(
    @{ Row = 2; Id = '#1.-borger'; Variant = '1. Borger' }
    @{ Row = 3; Id = '#haerlem'; Variant = 'Haerlem' }
    @{ Row = 4; Id = '#knecht'; Variant = 'Knecht' }
    @{ Row = 5; Id = '#g.-stuy'; Variant = 'G. Stuy' }
    @{ Row = 6; Id = '#2.-borg'; Variant = '2. Borg' }
    @{ Row = 7; Id = '#assend'; Variant = 'Assend' }
    @{ Row = 8; Id = '#2.-sold'; Variant = '2. Sold' }
    @{ Row = 9; Id = '#p.-kies'; Variant = 'P. Kies' }
    @{ Row = 10; Id = '#1.-vrouw'; Variant = '1. vrouw' }
    @{ Row = 11; Id = '#2.-vrouw'; Variant = '2. vrouw' }
    @{ Row = 12; Id = '#verlaen'; Variant = 'Verlaen' }
    @{ Row = 13; Id = '#de-maeght-haerlem'; Variant = 'De Maeght Haerlem' }
    @{ Row = 14; Id = '#aldeg'; Variant = 'Aldeg' }
    @{ Row = 15; Id = '#tijdt'; Variant = 'Tijdt' }
    @{ Row = 16; Id = '#maegera'; Variant = 'Maegera' }
    @{ Row = 17; Id = '#1.-borg'; Variant = '1. Borg' }
    @{ Row = 18; Id = '#vvaer-heyd'; Variant = 'VVaer-heyd' }
    @{ Row = 19; Id = '#steenba'; Variant = 'Steenba' }
    @{ Row = 20; Id = '#1.-vrouw'; Variant = '1. Vrouw' }
    @{ Row = 21; Id = '#rippard'; Variant = 'Rippard' }
    @{ Row = 22; Id = '#de-vries'; Variant = 'De Vries' }
    @{ Row = 23; Id = '#1.-sold'; Variant = '1. Sold' }
    @{ Row = 24; Id = '#2.-vrou'; Variant = '2. Vrou' }
    @{ Row = 25; Id = '#bordet'; Variant = 'Bordet' }
    @{ Row = 26; Id = '#prince-van-orangien'; Variant = 'Prince van Orangien' }
    @{ Row = 27; Id = '#m.-drag'; Variant = 'M. Drag' }
    @{ Row = 28; Id = '#1.-vrou'; Variant = '1. Vrou' }
    @{ Row = 29; Id = '#2.-borger'; Variant = '2. Borger' }
    @{ Row = 30; Id = '#i.-vliet'; Variant = 'I. Vliet' }
    @{ Row = 31; Id = '#alecto'; Variant = 'Alecto' }
    @{ Row = 32; Id = '#2.-vrouw'; Variant = '2. Vrouw' }
    @{ Row = 33; Id = '#i.-rome'; Variant = 'I. Rome' }
    @{ Row = 34; Id = '#schagen'; Variant = 'Schagen' }
    @{ Row = 35; Id = '#don-fr'; Variant = 'Don Fr' }
    @{ Row = 36; Id = '#noircar'; Variant = 'Noircar' }
)

foreach ($r in $rows) {
    $ws.Range("B$($r.Row)").Value = $r.Id
    $ws.Range("C$($r.Row)").Value = $r.Variant
}

# Rows 2-29 carried an "x" in the is_prefered column (D) in the old export;
# the new export leaves is_prefered blank for every row, so clear it.
$ws.Range("D2:D29").ClearContents()
